$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 654.8182
$ws.Range("I33").Value = 529.9259
$ws.Range("J33").Value = 1216.8334
$ws.Range("K33").Value = 529.9259
$ws.Range("L33").Value = 1216.8334
$ws.Range("M33").Value = -300.9259
$ws.Range("N33").Value = -1674.8334

$ws.Range("H76").Value = 7456.758
$ws.Range("I76").Value = 9742.4375
$ws.Range("J76").Value = 5305.5293
$ws.Range("K76").Value = 9742.4375
$ws.Range("L76").Value = 5305.5293
$ws.Range("M76").Value = -9427.4375
$ws.Range("N76").Value = -5935.5293

$ws.Range("H79").Value = 7456.758
$ws.Range("I79").Value = 9742.4375
$ws.Range("J79").Value = 5305.5293
$ws.Range("K79").Value = 9742.4375
$ws.Range("L79").Value = 5305.5293
$ws.Range("M79").Value = -8650.4375
$ws.Range("N79").Value = -7489.5293

$ws.Range("H107").Value = 407.84616
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 500.22223
$ws.Range("K107").Value = 200
$ws.Range("L107").Value = 500.22223
$ws.Range("M107").Value = 1720
$ws.Range("N107").Value = -4340.22223

$ws.Range("H116").Value = 4076.9092
$ws.Range("J116").Value = 3637.2307
$ws.Range("L116").Value = 3637.2307
$ws.Range("N116").Value = -10521.2307

$ws.Range("H137").Value = 3740.742
$ws.Range("I137").Value = 3503.32
$ws.Range("K137").Value = 10509.96
$ws.Range("M137").Value = -7959.960000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 611.0909
$ws.Range("I2").Value = 544.9524
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 544.9524
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -431.9524
$ws.Range("N2").Value = -2226

$ws.Range("H51").Value = 9750
$ws.Range("J51").Value = 9750
$ws.Range("L51").Value = 9750
$ws.Range("N51").Value = -11262

$ws.Range("H116").Value = 611.0909
$ws.Range("I116").Value = 544.9524
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 544.9524
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1749.0476
$ws.Range("N116").Value = -6588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 611.0909
$ws.Range("I3").Value = 544.9524
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 544.9524
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -430.9524
$ws.Range("N3").Value = -2228

$ws.Range("H86").Value = 8138.25
$ws.Range("I86").Value = 10234.333
$ws.Range("K86").Value = 10234.333
$ws.Range("M86").Value = -9111.333000000001

$ws.Range("H89").Value = 8138.25
$ws.Range("I89").Value = 10234.333
$ws.Range("K89").Value = 51171.665
$ws.Range("M89").Value = -45555.665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2732.3333
$ws.Range("I31").Value = 1553.4736
$ws.Range("K31").Value = 1553.4736
$ws.Range("M31").Value = -1258.4736

$ws.Range("H34").Value = 2732.3333
$ws.Range("I34").Value = 1553.4736
$ws.Range("K34").Value = 1553.4736
$ws.Range("M34").Value = -1351.4736

$ws.Range("H58").Value = 2287.2693
$ws.Range("I58").Value = 2575.7058
$ws.Range("K58").Value = 2575.7058
$ws.Range("M58").Value = -2372.7058

$ws.Range("H99").Value = 54170.58
$ws.Range("I99").Value = 101334.1
$ws.Range("J99").Value = 1766.6666
$ws.Range("K99").Value = 101334.1
$ws.Range("L99").Value = 1766.6666
$ws.Range("M99").Value = -99836.10000000001
$ws.Range("N99").Value = -4762.6666

$ws.Range("H107").Value = 282.15384
$ws.Range("I107").Value = 282.15384
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 282.15384
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1637.84616
$ws.Range("N107").ClearContents()

$ws.Range("H122").Value = 1626.2778
$ws.Range("I122").Value = 962.125
$ws.Range("J122").Value = 2157.6
$ws.Range("K122").Value = 2886.375
$ws.Range("L122").Value = 6472.799999999999
$ws.Range("M122").Value = -436.375
$ws.Range("N122").Value = -11372.8

$ws.Range("H126").Value = 54170.58
$ws.Range("I126").Value = 101334.1
$ws.Range("J126").Value = 1766.6666
$ws.Range("K126").Value = 304002.3
$ws.Range("L126").Value = 5299.9998
$ws.Range("M126").Value = -301532.3
$ws.Range("N126").Value = -10239.9998

$ws.Range("H132").Value = 1930
$ws.Range("J132").Value = 3461
$ws.Range("L132").Value = 10383
$ws.Range("N132").Value = -15443

$ws.Range("H134").Value = 1763.6984
$ws.Range("I134").Value = 1181.8529
$ws.Range("J134").Value = 2445.862
$ws.Range("K134").Value = 3545.5587
$ws.Range("L134").Value = 7337.586
$ws.Range("M134").Value = -1010.5587
$ws.Range("N134").Value = -12407.586

$ws.Range("H136").Value = 2287.2693
$ws.Range("I136").Value = 2575.7058
$ws.Range("K136").Value = 7727.117400000001
$ws.Range("M136").Value = -5177.117400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 14334
$ws.Range("I132").Value = 11002
$ws.Range("J132").Value = 16000
$ws.Range("K132").Value = 99018
$ws.Range("L132").Value = 144000
$ws.Range("M132").Value = -96488
$ws.Range("N132").Value = -149060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3666.6667
$ws.Range("I43").Value = 3666.6667
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3666.6667
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3515.6667
$ws.Range("N43").ClearContents()

$ws.Range("H70").Value = 4022.3518
$ws.Range("I70").Value = 4033.5278
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 4033.5278
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -3763.5278
$ws.Range("N70").Value = -4540

$ws.Range("H73").Value = 4022.3518
$ws.Range("I73").Value = 4033.5278
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 4033.5278
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -3097.5278
$ws.Range("N73").Value = -5872

$ws.Range("H132").Value = 3432.38
$ws.Range("I132").Value = 2989.5312
$ws.Range("K132").Value = 8968.5936
$ws.Range("M132").Value = -6438.5936

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 566
$ws.Range("I22").Value = 435.45456
$ws.Range("J22").Value = 2002
$ws.Range("K22").Value = 435.45456
$ws.Range("L22").Value = 2002
$ws.Range("M22").Value = -140.45456
$ws.Range("N22").Value = -2592

$ws.Range("H27").Value = 566
$ws.Range("I27").Value = 435.45456
$ws.Range("J27").Value = 2002
$ws.Range("K27").Value = 435.45456
$ws.Range("L27").Value = 2002
$ws.Range("M27").Value = -328.45456
$ws.Range("N27").Value = -2216

$ws.Range("H50").Value = 15004
$ws.Range("I50").Value = 9876
$ws.Range("J50").Value = 20132
$ws.Range("K50").Value = 9876
$ws.Range("L50").Value = 20132
$ws.Range("M50").Value = -9239
$ws.Range("N50").Value = -21406

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 14740
$ws.Range("J39").Value = 14740
$ws.Range("L39").Value = 14740
$ws.Range("N39").Value = -15566

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H43").Value = 13290
$ws.Range("J43").Value = 13290
$ws.Range("L43").Value = 13290
$ws.Range("N43").Value = -13588

$ws.Range("H82").Value = 25150.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 25150.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 25150.5
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -25916.5

$ws.Range("H85").Value = 25150.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 25150.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 25150.5
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -27802.5

$ws.Range("H132").Value = 1841.9344
$ws.Range("J132").Value = 2453
$ws.Range("L132").Value = 7359
$ws.Range("N132").Value = -12419

$ws.Range("H136").Value = 11234840
$ws.Range("I136").Value = 16146309
$ws.Range("J136").Value = 359443.4
$ws.Range("K136").Value = 48438927
$ws.Range("L136").Value = 1078330.2
$ws.Range("M136").Value = -48436377
$ws.Range("N136").Value = -1083430.2
